$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-89 down to 67-90
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record's data
$ws.Cells.Item(66, 1).Value = 11
$ws.Cells.Item(66, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(66, 3).Value = "Bíobío"
$ws.Cells.Item(66, 4).Value = 44951
$ws.Cells.Item(66, 5).Value = 8
$ws.Cells.Item(66, 6).Value = "Fruta"
$ws.Cells.Item(66, 7).Value = 100103
$ws.Cells.Item(66, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(66, 9).Value = 100103002
$ws.Cells.Item(66, 10).Value = "Ciruela"
$ws.Cells.Item(66, 11).Value = "Fortuna"
$ws.Cells.Item(66, 12).Value = "Primera"
$ws.Cells.Item(66, 13).Value = 170
$ws.Cells.Item(66, 14).Value = 10000
$ws.Cells.Item(66, 15).Value = 11000
$ws.Cells.Item(66, 16).Value = 10529
$ws.Cells.Item(66, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(66, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(66, 19).Value = 585
$ws.Cells.Item(66, 20).Value = 18
